$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

for ($r = 2; $r -le 120; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = -1 * $cell.Value2
}
